$wb = $excel.ActiveWorkbook

# Update the "Quirks" sheet fear labels to the new "Fear of X" wording
$ws = $wb.Worksheets.Item("Quirks")

$ws.Range("A3").Value = "Fear of the dark"
$ws.Range("A4").Value = "Fear of fire"
$ws.Range("A2").Value = "Fear of heights"
$ws.Range("A5").Value = "Fear of water/the ocean"
$ws.Range("A6").Value = "Fear of snakes"
$ws.Range("A7").Value = "Fear of spiders"
$ws.Range("A8").Value = "Fear of holes"

# Make Quirks the active sheet/tab and set its selection
$ws.Activate()
$ws.Range("A10").Select()
